# "Load rooming list v1"
#
# The Rooming sheet's header/binding row (row 4) is updated to reflect the
# new data-loading field names:
#   - A4: "Resource_id" -> "id"
#   - B4: (empty)        -> "Resource.Code"   (new binding column)
#   - S4: "Origin.Name"  -> "Country_origin.Name"
#
# Order matters here: Excel's shared-string table is rebuilt/compacted on
# save, dropping any strings that are no longer referenced ("Resource_id",
# "Origin.Name") and appending freshly-introduced strings in the order
# they were first assigned. To reproduce the expected order of new shared
# strings (Country_origin.Name before Resource.Code), set S4 before B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rooming")

$ws.Range("A4").Value = "id"
$ws.Range("S4").Value = "Country_origin.Name"
$ws.Range("B4").Value = "Resource.Code"
